$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{}
$rows[2] = @("FAPs", "Wnt7b", "Fzd10", "FAPs", 1, 0.3333333333333333, 0.072919, 0.218757, 0.1477750351608889, 0.1477750351608889, 1, 0.3333333333333333, 0.005445666666666667, 0.016337, 0.1242225162340131, 0.1242225162340131, 0.0003970925676666667, 0.003573833109, 0.01835698670425539, 0.01835698670425538)
$rows[3] = @("FAPs", "Wnt7b", "Fzd10", "Inflammatory-Mac", 1, 0.3333333333333333, 0.072919, 0.218757, 0.1477750351608889, 0.1477750351608889, 2, 0.6666666666666666, 0.03095266666666667, 0.092858, 0.706069315814286, 0.7060693158142859, 0.002257037500666666, 0.020313337506, 0.1043394179704809, 0.1043394179704809)
$rows[4] = @("FAPs", "Wnt7b", "Fzd10", "MuSCs", 1, 0.3333333333333333, 0.072919, 0.218757, 0.1477750351608889, 0.1477750351608889, 1, 0.3333333333333333, 0.007439666666666666, 0.022319, 0.169708167951701, 0.169708167951701, 0.0005424930536666666, 0.004882437483, 0.02507863048615265, 0.02507863048615265)
$rows[5] = @("MuSCs", "Wnt7b", "Fzd10", "FAPs", 3, 1, 0.420527, 1.261581, 0.852224964839111, 0.852224964839111, 1, 0.3333333333333333, 0.005445666666666667, 0.016337, 0.1242225162340131, 0.1242225162340131, 0.002290049866333334, 0.020610448797, 0.1058655295297577, 0.1058655295297577)
$rows[6] = @("MuSCs", "Wnt7b", "Fzd10", "Inflammatory-Mac", 3, 1, 0.420527, 1.261581, 0.852224964839111, 0.852224964839111, 2, 0.6666666666666666, 0.03095266666666667, 0.092858, 0.706069315814286, 0.7060693158142859, 0.01301643205533333, 0.117147888498, 0.601729897843805, 0.6017298978438049)
$rows[7] = @("MuSCs", "Wnt7b", "Fzd10", "MuSCs", 3, 1, 0.420527, 1.261581, 0.852224964839111, 0.852224964839111, 1, 0.3333333333333333, 0.007439666666666666, 0.022319, 0.169708167951701, 0.169708167951701, 0.003128580704333334, 0.028157226339, 0.1446295374655483, 0.1446295374655483)

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $arr = New-Object "object[,]" 1,20
    for ($i = 0; $i -lt 20; $i++) {
        $arr[0,$i] = $vals[$i]
    }
    $rangeAddr = "A" + $r + ":T" + $r
    $ws.Range($rangeAddr).Value = $arr
}
